$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H116").Value = 5
